# Update countries & provincias Spain
#
# This refreshes the COVID-19 "Pais" dashboard data: updated case counts for a
# handful of countries (which also changes the descending sort order of the
# table by "Casos totales", causing a few rows to swap contents), plus the
# "last updated" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last refreshed" timestamp cell (A1) ---
# Row 1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 15 de Abril de 2020 a las 23:52"

# Row 4
$ws.Cells.Item(4, 2).Value = 641397
$ws.Cells.Item(4, 3).Value = 27511
$ws.Cells.Item(4, 4).Value = 48105
$ws.Cells.Item(4, 5).Value = 564898
$ws.Cells.Item(4, 7).Value = 2347
$ws.Cells.Item(4, 8).Value = 28394

# Row 95
$ws.Cells.Item(95, 2).Value = 542
$ws.Cells.Item(95, 3).Value = 14
$ws.Cells.Item(95, 4).Value = 226
$ws.Cells.Item(95, 5).Value = 284
$ws.Cells.Item(95, 7).Value = 2
$ws.Cells.Item(95, 8).Value = 32

# Row 131
$ws.Cells.Item(131, 1).Value = "Jamaica"
$ws.Cells.Item(131, 2).Value = 125
$ws.Cells.Item(131, 3).Value = 20
$ws.Cells.Item(131, 4).Value = 21
$ws.Cells.Item(131, 5).Value = 99
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = 5

# Row 132
$ws.Cells.Item(132, 1).Value = "Camboya"
$ws.Cells.Item(132, 2).Value = 122
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 96
$ws.Cells.Item(132, 5).Value = 26
$ws.Cells.Item(132, 6).Value = 1
$ws.Cells.Item(132, 8).Value = 0

# Row 133
$ws.Cells.Item(133, 1).Value = "Congo"
$ws.Cells.Item(133, 2).Value = 117
$ws.Cells.Item(133, 3).Value = 43
$ws.Cells.Item(133, 4).Value = 11
$ws.Cells.Item(133, 5).Value = 101
$ws.Cells.Item(133, 8).Value = 5

# Row 134
$ws.Cells.Item(134, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(134, 2).Value = 114
$ws.Cells.Item(134, 3).Value = 1
$ws.Cells.Item(134, 4).Value = 19
$ws.Cells.Item(134, 5).Value = 87
$ws.Cells.Item(134, 6).Value = 0
$ws.Cells.Item(134, 8).Value = 8

# Row 135
$ws.Cells.Item(135, 1).Value = "Madagascar"
$ws.Cells.Item(135, 2).Value = 110
$ws.Cells.Item(135, 3).Value = 2
$ws.Cells.Item(135, 4).Value = 29
$ws.Cells.Item(135, 5).Value = 81
$ws.Cells.Item(135, 6).Value = 1
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 0

# Row 141
$ws.Cells.Item(141, 1).Value = "Bermudas"
$ws.Cells.Item(141, 3).Value = 24
$ws.Cells.Item(141, 4).Value = 33
$ws.Cells.Item(141, 6).Value = 3
$ws.Cells.Item(141, 8).Value = 5

# Row 142
$ws.Cells.Item(142, 1).Value = "Togo"
$ws.Cells.Item(142, 2).Value = 81
$ws.Cells.Item(142, 3).Value = 4
$ws.Cells.Item(142, 4).Value = 35
$ws.Cells.Item(142, 5).Value = 43
$ws.Cells.Item(142, 8).Value = 3

# Row 143
$ws.Cells.Item(143, 1).Value = "Gabon"
$ws.Cells.Item(143, 3).Value = 23
$ws.Cells.Item(143, 4).Value = 4
$ws.Cells.Item(143, 5).Value = 75
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 1

# Row 144
$ws.Cells.Item(144, 1).Value = "Somalia"
$ws.Cells.Item(144, 2).Value = 80
$ws.Cells.Item(144, 3).Value = 20
$ws.Cells.Item(144, 4).Value = 2
$ws.Cells.Item(144, 5).Value = 73
$ws.Cells.Item(144, 6).Value = 2
$ws.Cells.Item(144, 7).Value = 3
$ws.Cells.Item(144, 8).Value = 5

# Row 145
$ws.Cells.Item(145, 1).Value = "Liechtenstein"
$ws.Cells.Item(145, 2).Value = 79
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 55
$ws.Cells.Item(145, 5).Value = 23
$ws.Cells.Item(145, 8).Value = 1

# Row 146
$ws.Cells.Item(146, 1).Value = "Birmania"
$ws.Cells.Item(146, 2).Value = 74
$ws.Cells.Item(146, 3).Value = 11
$ws.Cells.Item(146, 4).Value = 2
$ws.Cells.Item(146, 5).Value = 68
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 8).Value = 4

# Row 147
$ws.Cells.Item(147, 1).Value = "Barbados"
$ws.Cells.Item(147, 2).Value = 73
$ws.Cells.Item(147, 4).Value = 15
$ws.Cells.Item(147, 5).Value = 53
$ws.Cells.Item(147, 6).Value = 4
$ws.Cells.Item(147, 8).Value = 5

# Row 148
$ws.Cells.Item(148, 1).Value = "Liberia"
$ws.Cells.Item(148, 2).Value = 59
$ws.Cells.Item(148, 4).Value = 4
$ws.Cells.Item(148, 5).Value = 49
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 8).Value = 6

# Row 172
$ws.Cells.Item(172, 1).Value = "Zimbabue"
$ws.Cells.Item(172, 3).Value = 6
$ws.Cells.Item(172, 4).Value = 1
$ws.Cells.Item(172, 5).Value = 19
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 8).Value = 3

# Row 173
$ws.Cells.Item(173, 1).Value = "Antigua y Barbuda"
$ws.Cells.Item(173, 2).Value = 23
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 3
$ws.Cells.Item(173, 5).Value = 18
$ws.Cells.Item(173, 6).Value = 1
$ws.Cells.Item(173, 8).Value = 2

# Row 174
$ws.Cells.Item(174, 1).Value = "Maldivas"
$ws.Cells.Item(174, 2).Value = 22
$ws.Cells.Item(174, 3).Value = 2
$ws.Cells.Item(174, 4).Value = 16
$ws.Cells.Item(174, 5).Value = 6

# Row 175
$ws.Cells.Item(175, 1).Value = "Laos"
$ws.Cells.Item(175, 4).Value = 1
$ws.Cells.Item(175, 5).Value = 18
$ws.Cells.Item(175, 8).Value = 0

# Row 176
$ws.Cells.Item(176, 1).Value = "Angola"
$ws.Cells.Item(176, 2).Value = 19
$ws.Cells.Item(176, 4).Value = 5
$ws.Cells.Item(176, 5).Value = 12
$ws.Cells.Item(176, 8).Value = 2

# Row 177
$ws.Cells.Item(177, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(177, 4).Value = 1
$ws.Cells.Item(177, 5).Value = 17
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 8).Value = 0

# Row 178
$ws.Cells.Item(178, 1).Value = "Belice"
$ws.Cells.Item(178, 3).Value = 0
$ws.Cells.Item(178, 4).Value = 0
$ws.Cells.Item(178, 5).Value = 16
$ws.Cells.Item(178, 6).Value = 1
$ws.Cells.Item(178, 8).Value = 2

# Row 187
$ws.Cells.Item(187, 1).Value = "Granada"
$ws.Cells.Item(187, 6).Value = 2

# Row 188
$ws.Cells.Item(188, 1).Value = "San Cristobal y Nieves"
$ws.Cells.Item(188, 6).Value = 0

# Row 195
$ws.Cells.Item(195, 1).Value = "Montserrat"
$ws.Cells.Item(195, 6).Value = 1

# Row 196
$ws.Cells.Item(196, 1).Value = "Islas Malvinas"
$ws.Cells.Item(196, 6).Value = 0
